# Add a new data row (row 7) to the "Artfynd" sheet, mirroring the
# structure of the existing rows (e.g. row 6), with new observation data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric columns
$ws.Range("A7").Value = 131257158
$ws.Range("B7").Value = 57881
$ws.Range("E7").Value = 100049
$ws.Range("Q7").Value = 567581
$ws.Range("R7").Value = 6509575
$ws.Range("S7").Value = 10

# Text columns
$ws.Range("D7").Value = "NT"
$ws.Range("F7").Value = "Spillkråka"
$ws.Range("G7").Value = "Dryocopus martius"
$ws.Range("H7").Value = "(Linnaeus, 1758)"
$ws.Range("M7").Value = "färska spår"
$ws.Range("P7").Value = "Sjöberga 1:2, Ög"
$ws.Range("T7").Value = "Östergötland"
$ws.Range("U7").Value = "Norrköping"
$ws.Range("V7").Value = "Östergötland"
$ws.Range("W7").Value = "Kvillinge"
$ws.Range("AW7").Value = "Anette Källman"
$ws.Range("AX7").Value = "Anette Källman"

# Date columns stored as plain text (leading apostrophe forces text,
# matching the source data which keeps these as literal strings rather
# than date serials)
$ws.Range("Y7").Value = "'2026-02-21"
$ws.Range("AA7").Value = "'2026-02-21"

# Boolean columns
$ws.Range("AD7").Value = $false
$ws.Range("AE7").Value = $false
$ws.Range("AG7").Value = $false

# Columns that are present on the row but hold an empty string (as on the
# neighbouring rows) rather than being entirely absent.
$ws.Range("I7").Formula = "="""""
$ws.Range("K7").Formula = "="""""
$ws.Range("L7").Formula = "="""""
$ws.Range("N7").Formula = "="""""
$ws.Range("AT7").Formula = "="""""
$ws.Range("AY7").Formula = "="""""
